$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 49 (columns L-P) ---
$ws.Range("L49").Value = '32 batch: 0.9611137218045113, 2 batch: 0.9399906781635983'
$ws.Range("M49").Value = '32 batch: 0.9464285714285714, 2 batch: 0.958168259147052'
$ws.Range("N49").Value = '32 batch: 0.08223684210526316, 2 batch: 0.6313213703099511'
$ws.Range("O49").Value = '32 batch: 0.08634868421052631, 2 batch: 0.636215334420881'
$ws.Range("P49").Value = 'slurm-42391225'

# --- Row 50 ---
$ws.Range("A50").Value = '1-7 (training), 0 (testing)'
$ws.Range("B50").Value = '16 words'
$ws.Range("C50").Value = '(detrended) gaussian weighted 3D fMRI image 2-8 seconds after each word'
$ws.Range("D50").Value = 'Default CLIP'
$ws.Range("E50").Value = 'Default CLIP'
$ws.Range("F50").Value = '3D Resnet18'
$ws.Range("G50").Value = 'Cosine Similarity'
$ws.Range("H50").Value = 'embed_dim=512, image_resolution, layers=(2,2,2,2), width=64, context_length=16, vocab_size, transformer_width, transformer_heads, transformer_layers'
$ws.Range("I50").Value = 'LR=5e-4, batch_size=32, weight_decay=0.2'
$ws.Range("K50").Value = '700/0/100'
$ws.Range("L50").Value = '32 batch: 0.03125, 2 batch: 0.5'
$ws.Range("M50").Value = '32 batch: 0.03125, 2 batch: 0.5'
$ws.Range("N50").Value = '32 batch: 0.03125, 2 batch: 0.5'
$ws.Range("O50").Value = '32 batch: 0.03125, 2 batch: 0.5'
$ws.Range("P50").Value = 'slurm-42396686'
$ws.Range("J50").Value = 100

# --- Row 51 ---
$ws.Range("A51").Value = '1-7 (training), 0 (testing)'
$ws.Range("B51").Value = '16 words'
$ws.Range("C51").Value = '(detrended) gaussian weighted 3D fMRI image 2-8 seconds after each word'
$ws.Range("D51").Value = 'Default CLIP'
$ws.Range("E51").Value = 'Default CLIP'
$ws.Range("F51").Value = '3D Resnet18'
$ws.Range("G51").Value = 'Cosine Similarity'
$ws.Range("H51").Value = 'embed_dim=512, image_resolution, layers=(2,2,2,2), width=64, context_length=16, vocab_size, transformer_width, transformer_heads, transformer_layers'
$ws.Range("I51").Value = 'LR=1e-5, batch_size=32, weight_decay=0.2'
$ws.Range("K51").Value = '700/0/100'
$ws.Range("L51").Value = '32 batch: 0.9459586466165414, 2 batch: 0.9461663947797716'
$ws.Range("M51").Value = '32 batch: 0.930921052631579, 2 batch: 0.9550221393614542'
$ws.Range("N51").Value = '32 batch: 0.07072368421052631, 2 batch: 0.5628058727569332'
$ws.Range("O51").Value = '32 batch: 0.0625, 2 batch: 0.5497553017944535'
$ws.Range("P51").Value = 'slurm-42396686'
$ws.Range("J51").Value = 100

# --- Row 52 ---
$ws.Range("A52").Value = '1-7 (training), 0 (testing)'
$ws.Range("B52").Value = '16 words'
$ws.Range("C52").Value = '(detrended) gaussian weighted 3D fMRI image 2-8 seconds after each word'
$ws.Range("D52").Value = 'Default CLIP'
$ws.Range("E52").Value = 'Default CLIP'
$ws.Range("F52").Value = '3D Resnet18'
$ws.Range("G52").Value = 'Cosine Similarity'
$ws.Range("H52").Value = 'embed_dim=512, image_resolution, layers=(2,2,2,2), width=64, context_length=16, vocab_size, transformer_width, transformer_heads, transformer_layers'
$ws.Range("I52").Value = 'LR=5e-5, batch_size=32, weight_decay=0.2'
$ws.Range("K52").Value = '700/0/100'
$ws.Range("L52").Value = '32 batch: 0.7095864661654135, 2 batch: 0.9591004427872291'
$ws.Range("M52").Value = '32 batch: 0.6766917293233082, 2 batch: 0.9882311815427639'
$ws.Range("N52").Value = '32 batch: 0.10032894736842106, 2 batch: 0.6419249592169658'
$ws.Range("O52").Value = '32 batch: 0.09046052631578948, 2 batch'
$ws.Range("P52").Value = 'slurm-42396686'
$ws.Range("J52").Value = 100

# --- Row 53 ---
$ws.Range("A53").Value = '1-7 (training), 0 (testing)'
$ws.Range("B53").Value = '16 words'
$ws.Range("C53").Value = '(detrended) gaussian weighted 3D fMRI image 2-8 seconds after each word'
$ws.Range("D53").Value = 'Default CLIP'
$ws.Range("E53").Value = 'Default CLIP'
$ws.Range("F53").Value = '3D Resnet18'
$ws.Range("G53").Value = 'Cosine Similarity'
$ws.Range("H53").Value = 'embed_dim=1024, image_resolution, layers=(2,2,2,2), width=64, context_length=16, vocab_size, transformer_width, transformer_heads, transformer_layers'
$ws.Range("I53").Value = 'LR=1e-5, batch_size=64, weight_decay=0.2'
$ws.Range("K53").Value = '700/0/100'
$ws.Range("L53").Value = '64 batch: 0.945371240601503, 2 batch: 0.9108599394080634'
$ws.Range("M53").Value = '64 batch: 0.9273966165413534, 2 batch: 0.9535073409461664'
$ws.Range("N53").Value = '64 batch: 0.017269736842105265, 2 batch: 0.586460032626427'
$ws.Range("O53").Value = '64 batch: 0.02631578947368421, 2 batch: 0.6117455138662317'
$ws.Range("P53").Value = 'slurm-42435618'
$ws.Range("J53").Value = 100

# --- Row 54 ---
$ws.Range("A54").Value = '1-7 (training), 0 (testing)'
$ws.Range("B54").Value = '16 words'
$ws.Range("C54").Value = '(detrended) gaussian weighted 3D fMRI image 2-8 seconds after each word'
$ws.Range("D54").Value = 'Default CLIP'
$ws.Range("E54").Value = 'Default CLIP'
$ws.Range("F54").Value = '3D Resnet18'
$ws.Range("G54").Value = 'Cosine Similarity'
$ws.Range("H54").Value = 'embed_dim=1024, image_resolution, layers=(2,2,2,2), width=64, context_length=16, vocab_size, transformer_width, transformer_heads, transformer_layers'
$ws.Range("I54").Value = 'LR=5e-5, batch_size=64, weight_decay=0.2'
$ws.Range("K54").Value = '700/0/100'
$ws.Range("L54").Value = '64 batch: 0.543233082706767, 2 batch: 0.9364949895129341'
$ws.Range("M54").Value = '64 batch: 0.5110432330827067, 2 batch: 0.97657888604055'
$ws.Range("N54").Value = '64 batch: 0.05674342105263158, 2 batch: 0.6786296900489397'
$ws.Range("O54").Value = '64 batch: 0.0805921052631579, 2 batch: 0.7047308319738989'
$ws.Range("P54").Value = 'slurm-42435618'
$ws.Range("J54").Value = 100

# --- Row 55 ---
$ws.Range("C55").Value = '(fmri channel for each word) (detrended) gaussian weighted 3D fMRI image 2-8 seconds after each word'

# --- View/selection updates ---
$ws.Activate() | Out-Null
$av = $excel.ActiveWindow
$av.ScrollRow = 45
$av.ScrollColumn = 1
$ws.Range("C50").Select() | Out-Null
